$d = $word.ActiveDocument

# wdYellow highlight color index
$wdYellow = 7

function Set-ParaText($para, [string]$text, [bool]$highlight) {
    $r = $para.Range
    $r.End = $r.End - 1
    $r.Text = $text
    if ($highlight) {
        $r2 = $para.Range
        $r2.End = $r2.End - 1
        $r2.HighlightColorIndex = $wdYellow
    }
}

# Process from the LAST paragraph to the FIRST so that deleting the
# blank separator paragraphs does not shift the index of paragraphs
# that have not been processed yet.

# Paragraph 13: "[6] I. Hallmann ..." -> "[8][6] I. Hallmann ..."
$p13 = $d.Paragraphs.Item(13)
Set-ParaText $p13 '[8][6] I. Hallmann and B. Siemiatkowska, \u201cArtificial landmark navigation system,\u201d in Proc. Int. Symp. Intell. Robot. Syst., Jul. 2001, pp. 219\u2013228.' $false

# Paragraph 12: empty separator -> delete
$d.Paragraphs.Item(12).Range.Delete()

# Paragraph 11: "[5] Gueaieb ..." -> "[7][5] Gueaieb ..."
$p11 = $d.Paragraphs.Item(11)
Set-ParaText $p11 '[7][5] Gueaieb, W.; Miah, Md.S., "An Intelligent Mobile Robot Navigation Technique Using RFID Technology," Instrumentation and Measurement, IEEE Transactions on , vol.57, no.9, pp.1908,1917, Sept. 2008' $false

# Paragraph 10: empty separator -> delete
$d.Paragraphs.Item(10).Range.Delete()

# Paragraph 9: "[4] Hahnel ..." -> "[6][4] Hahnel ..." (also drops the
# stray space before the comma after "on" and adds a space in "pp.1015, 1020")
$p9 = $d.Paragraphs.Item(9)
Set-ParaText $p9 '[6][4] Hahnel, D.; Burgard, W.; Fox, D.; Fishkin, K.; Philipose, M., "Mapping and localization with RFID technology," Robotics and Automation, 2004. Proceedings. ICRA ''04. 2004 IEEE International Conference on, vol.1, no., pp.1015, 1020 Vol.1, 26 April-1 May 2004' $false

# Paragraph 8: empty separator -> delete
$d.Paragraphs.Item(8).Range.Delete()

# Paragraph 7: "[3] Klaus Finkenzeller ..." -> "[5] Klaus Finkenzeller ..."
$p7 = $d.Paragraphs.Item(7)
Set-ParaText $p7 '[5] Klaus Finkenzeller. RFID Handboook: Radio-Frequency Identification Fundamentals and Applications. Wiley, New York, 2000.' $false

# Paragraph 6: empty separator -> delete
$d.Paragraphs.Item(6).Range.Delete()

# Paragraph 5: "[3] Tripathy ..." -> brand-new reference "[4] Want, R. ..." highlighted yellow
$p5 = $d.Paragraphs.Item(5)
Set-ParaText $p5 '[4] Want, R., "An introduction to RFID technology," Pervasive Computing, IEEE , vol.5, no.1, pp.25,33, Jan.-March 2006' $true

# Paragraph 4: empty separator -> delete
$d.Paragraphs.Item(4).Range.Delete()

# Paragraph 3: "[2] Gijeong ..." -> "[7] Tripathy ..." highlighted yellow
$p3 = $d.Paragraphs.Item(3)
Set-ParaText $p3 '[7] Tripathy, H. K., Tripathy, B. K., & Das, P. K. (2008). A Prospective Fuzzy Logic approach to Knowledge-based Navigation of Mobile LEGO-Robot. Journal of Convergence Information Technology, 3(1), 64-70.' $true

# Paragraph 2: empty -> fill with the (trimmed) former [2] Gijeong entry
$p2 = $d.Paragraphs.Item(2)
Set-ParaText $p2 '[2] Gijeong Jang; Sungho Lee; Inso Kweon, "Color landmark based self-localization for indoor mobile robots," Robotics and Automation, 2002. Proceedings. ICRA ''02. IEEE International Conference on, vol.1, no.' $false

# Paragraph 1: merge the runs of the [1] Azlan reference into one run
$p1 = $d.Paragraphs.Item(1)
Set-ParaText $p1 ' [1] Azlan, N.Z.; Zainudin, F.; Yusuf, H.M.; Toha, S.F.; Yusoff, S.Z.S.; Osman, N.H., "Fuzzy Logic Controlled Miniature LEGO Robot for Undergraduate Training System," Industrial Electronics and Applications, 2007. ICIEA 2007. 2nd IEEE Conference on, vol., no., pp.2184, 2188, 23-25 May 2007' $false

Write-Output ("Paragraph count: " + $d.Paragraphs.Count)
